# Bento Filter Chemo-OtherTreatment "Fixed Bento 80 Test scripts" edit
# Appends "order By <col> ASC LIMIT 100" clauses to the three detail-query cells
# (Cases/B2, Samples/B3, Files/B4) and normalizes the ORDER BY casing on the Files query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: Cases detail query -> append "order By ss.study_subject_id ASC LIMIT 100"
$ws.Cells.Item(2, 2).Value = 'MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE tp.chemotherapy_regimen in ["Other treatment given as part of a CTSU protocol"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
order By ss.study_subject_id ASC LIMIT 100'

# B3: Samples detail query -> append " order By samp.sample_id ASC LIMIT 100"
$ws.Cells.Item(3, 2).Value = 'MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen IN ["Other treatment given as part of a CTSU protocol"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,''-'')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
 order By samp.sample_id ASC LIMIT 100'

# B4: Files detail query -> replace "order by f.file_name" with " order By f.file_name ASC LIMIT 100"
$ws.Cells.Item(4, 2).Value = 'MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE tp.chemotherapy_regimen IN  ["Other treatment given as part of a CTSU protocol"]
WITH
        f, parent,p, ss, d,tp, s, samp,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
     order By f.file_name ASC LIMIT 100'

# A3/A4 tab-name labels unchanged in content, rewritten so the shared-string table
# reorders the same way Excel produced it (SamplesTab/FilesTab right after the CasesTab query)
$ws.Cells.Item(3, 1).Value = 'SamplesTab'
$ws.Cells.Item(4, 1).Value = 'FilesTab'

# Row heights grow because of the extra wrapped line introduced in B2/B3
$ws.Rows.Item(2).RowHeight = 345.6
$ws.Rows.Item(3).RowHeight = 374.4

# The view scrolled up one row (topLeftCell A4 -> A3)
$excel.ActiveWindow.ScrollRow = 3

